$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "'001"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 796946022.88
$ws.Range("P2").Value = 139971788.11
$ws.Range("Q2").Value = 113744021.01
$ws.Range("R2").Value = 928.9079349734
$ws.Range("S2").Value = 153476622.85
$ws.Range("T2").Value = 19.9369371477
$ws.Range("U2").Value = 152606762.75
$ws.Range("V2").Value = 9.315527354
$ws.Range("W2").Value = 352495723.78
$ws.Range("X2").Value = 141697375.41
$ws.Range("Y2").Value = 1.2795549971
$ws.Range("Z2").Value = 5730141.43
$ws.Range("AA2").Value = -20.1217830597
$ws.Range("AB2").Value = 444450299.1
$ws.Range("AC2").Value = 66.9422343972
$ws.Range("AD2").Value = 50.6118444849
$ws.Range("AE2").Value = 34.075174711
$ws.Range("AF2").Value = 151.0069815155
$ws.Range("AG2").Value = 44.2308153451
